# "Update data" - append the 2020-04-21 row to the Confirmados and Mortes
# sheets (the Populacao sheet's UF/Populacao header values are unaffected
# in content - only their underlying shared-string index shifts, which
# Excel manages automatically via the shared string table).

$wb = $excel.ActiveWorkbook
$wsConfirmados = $wb.Worksheets.Item("Confirmados")
$wsMortes = $wb.Worksheets.Item("Mortes")

$newDate = "2020-04-21"
$newRow = 58

# New row of confirmed-case counts, one value per state column (B..AB),
# in the same column order as row 1's header.
$confirmados = @(195,210,457,2270,1489,3716,881,1212,421,1396,181,173,1230,1026,263,1024,2908,186,5306,608,904,199,247,1063,15385,92,37)

# New row of death counts, one value per state column (B..AB).
$mortes = @(8,19,13,193,47,215,24,34,19,60,6,6,44,38,33,51,260,14,461,28,27,4,3,35,1093,5,1)

# Write the date as text (not an auto-converted date serial) by briefly
# forcing a text number format, then restoring the default "Normal"
# style so the cell matches the plain t="s" cells used elsewhere in
# column A.
$wsConfirmados.Cells.Item($newRow, 1).NumberFormat = "@"
$wsConfirmados.Cells.Item($newRow, 1).Value = $newDate
$wsConfirmados.Cells.Item($newRow, 1).Style = "Normal"
for ($i = 0; $i -lt $confirmados.Length; $i++) {
    $wsConfirmados.Cells.Item($newRow, $i + 2).Value = $confirmados[$i]
}

$wsMortes.Cells.Item($newRow, 1).NumberFormat = "@"
$wsMortes.Cells.Item($newRow, 1).Value = $newDate
$wsMortes.Cells.Item($newRow, 1).Style = "Normal"
for ($i = 0; $i -lt $mortes.Length; $i++) {
    $wsMortes.Cells.Item($newRow, $i + 2).Value = $mortes[$i]
}
